$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.062.65"
$ws.Range("E2").Value = "  +3.05%  "

$ws.Range("D3").Value = "1.596.62"
$ws.Range("E3").Value = "  +1.89%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'212.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.21%  "

$ws.Range("D7").Value = "'0.484"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.30%  "

$ws.Range("E8").Value = "  +2.35%  "

$ws.Range("D9").Value = "'0.0615"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.22%  "

$ws.Range("D10").Value = "'17.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.52%  "

$ws.Range("E11").Value = "  +4.36%  "

$ws.Range("D12").Value = "1.818.60"
$ws.Range("E12").Value = "  +1.91%  "

$ws.Range("D13").Value = "1.587.31"
$ws.Range("E13").Value = "  +1.22%  "

$ws.Range("E14").Value = "  -0.58%  "

$ws.Range("D15").Value = "'0.509"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.59%  "

$ws.Range("D16").Value = "26.037.03"
$ws.Range("E16").Value = "  +2.94%  "

$ws.Range("E18").Value = "  +1.05%  "

$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("D20").Value = "'203.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.68%  "

$ws.Range("E21").Value = "  +2.24%  "

$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("E23").Value = "  +1.47%  "

$ws.Range("D25").Value = "'141.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.84%  "

$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("D27").Value = "'0.124"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.14%  "

$ws.Range("D28").Value = "'15.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.41%  "

$ws.Range("D29").Value = "'6.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.40%  "

$ws.Range("E30").Value = "  +1.29%  "

$ws.Range("D31").Value = "'0.0469"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.30%  "

$ws.Range("E32").Value = "  +2.75%  "

$ws.Range("D33").Value = "'2.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.66%  "

$ws.Range("E34").Value = "  +0.80%  "

$ws.Range("D35").Value = "'2.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.72%  "

$ws.Range("D36").Value = "1.108.82"
$ws.Range("E36").Value = "  +1.84%  "

$ws.Range("E37").Value = "  +6.87%  "

$ws.Range("E38").Value = "  +0.29%  "

$ws.Range("E39").Value = "  +0.52%  "

$ws.Range("E40").Value = "  +0.75%  "

$ws.Range("D41").Value = "'0.491"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.87%  "

$ws.Range("D42").Value = "'0.776"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.48%  "

$ws.Range("D43").Value = "1.730.59"
$ws.Range("E43").Value = "  +1.87%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.36%  "

$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'92.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.17%  "

$ws.Range("D46").Value = "'1.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.41%  "

$ws.Range("D47").Value = "'53.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.10%  "

$ws.Range("D48").Value = "0.0₇0989"
$ws.Range("E48").Value = "  -11.34%  "

$ws.Range("D49").Value = "'0.0504"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.05%  "

$ws.Range("E50").Value = "  +0.71%  "

$ws.Range("E51").Value = "  +0.14%  "
